$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks up front to avoid orphaned refs after row deletion
$ws.Hyperlinks.Delete()

# --- Update data rows 2-9 with new content ---
# Row 2
$ws.Range('A2').Value = '2025-10-25 06:24:40'
$ws.Range('B2').Value = '【時給1,600円 / 学生限定】AIでプロダクトを生成したことがある学生の方を大募集!!'
$ws.Range('C2').Value = 'システム開発'
$ws.Range('D2').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E2').Value = '期限情報なし'
$ws.Range('F2').Value = 'https://www.lancers.jp/work/detail/5420120'
$ws.Range('G2').Value = 303
$ws.Range('H2').Value = '🔥AI,Ai'

# Row 3
$ws.Range('A3').Value = '2025-10-25 06:24:40'
$ws.Range('B3').Value = '【学生発スタートアップ】留学×住まいマッチングアプリ開発仲間募集'
$ws.Range('C3').Value = 'システム開発'
$ws.Range('D3').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E3').Value = '期限情報なし'
$ws.Range('F3').Value = 'https://www.lancers.jp/work/detail/5420198'
$ws.Range('G3').Value = 100
$ws.Range('H3').Value = '◆開発 ◇アプリ'

# Row 4
$ws.Range('A4').Value = '2025-10-25 06:24:40'
$ws.Range('B4').Value = '【クリニック向け】セキュアなPC管理の遠隔保守方法を教えてください(助言のみでのお支払い)'
$ws.Range('C4').Value = 'システム開発'
$ws.Range('D4').Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range('E4').Value = '期限情報なし'
$ws.Range('F4').Value = 'https://www.lancers.jp/work/detail/5420306'
$ws.Range('G4').Value = 30
$ws.Range('H4').Value = '◇管理'

# Row 5
$ws.Range('A5').Value = '2025-10-25 06:24:40'
$ws.Range('B5').Value = '【Ubuntu】MySQLデータを自動CSV化しクラウド保存構築'
$ws.Range('C5').Value = 'システム開発'
$ws.Range('D5').Value = '~ 5,000 円 / 固定'
$ws.Range('E5').Value = '期限情報なし'
$ws.Range('F5').Value = 'https://www.lancers.jp/work/detail/5420180'
$ws.Range('G5').Value = 30
$ws.Range('H5').Value = '◇MySQL'

# Row 6
$ws.Range('A6').Value = '2025-10-25 06:24:40'
$ws.Range('B6').Value = '【急募】時間単位で入札できるシステム構築の依頼'
$ws.Range('C6').Value = 'システム開発'
$ws.Range('D6').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E6').Value = '期限情報なし'
$ws.Range('F6').Value = 'https://www.lancers.jp/work/detail/5411365'
$ws.Range('G6').Value = 40
$ws.Range('H6').ClearContents()

# Row 7
$ws.Range('A7').Value = '2025-10-25 06:24:40'
$ws.Range('B7').Value = '運用中HPのドメイン分け'
$ws.Range('C7').Value = 'システム開発'
$ws.Range('D7').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E7').Value = '期限情報なし'
$ws.Range('F7').Value = 'https://www.lancers.jp/work/detail/5420277'
$ws.Range('G7').Value = 13
$ws.Range('H7').ClearContents()

# Row 8
$ws.Range('A8').Value = '2025-10-25 06:24:40'
$ws.Range('B8').Value = '注目 【急募】YouTubeの音楽配信構築の依頼です'
$ws.Range('C8').Value = 'システム開発'
$ws.Range('D8').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E8').Value = '期限情報なし'
$ws.Range('F8').Value = 'https://www.lancers.jp/work/detail/5420233'
$ws.Range('G8').Value = 13
$ws.Range('H8').ClearContents()

# Row 9
$ws.Range('A9').Value = '2025-10-25 06:24:40'
$ws.Range('B9').Value = '【急募】イベント用問い合わせLINE構築のフリーランス募集!'
$ws.Range('C9').Value = 'システム開発'
$ws.Range('D9').Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range('E9').Value = '期限情報なし'
$ws.Range('F9').Value = 'https://www.lancers.jp/work/detail/5420186'
$ws.Range('G9').Value = 10
$ws.Range('H9').ClearContents()

# --- Delete now-obsolete rows 10:19 ---
$ws.Range("A10:H19").EntireRow.Delete()

# --- Re-create hyperlinks for F2:F9 with final URLs ---
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://www.lancers.jp/work/detail/5420120') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://www.lancers.jp/work/detail/5420198') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://www.lancers.jp/work/detail/5420306') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F5'), 'https://www.lancers.jp/work/detail/5420180') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F6'), 'https://www.lancers.jp/work/detail/5411365') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F7'), 'https://www.lancers.jp/work/detail/5420277') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F8'), 'https://www.lancers.jp/work/detail/5420233') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F9'), 'https://www.lancers.jp/work/detail/5420186') | Out-Null

# --- Column width adjustments (raw xlsx width = ColumnWidth + 0.83) ---
$ws.Columns.Item(2).ColumnWidth = 46.17   # B: 50 -> 47
$ws.Columns.Item(4).ColumnWidth = 29.17   # D: 41 -> 30
$ws.Columns.Item(8).ColumnWidth = 11.17   # H: 18 -> 12

